$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$shp = $s.Shapes.Item(3)

# --- Resize / reposition the "first element" label textbox ---
# Old:  off  x=2486189 y=5578351   ext cx=5729958 cy=927101
# New:  off  x=2486189 y=5575107   ext cx=4408258 cy=933589
$shp.Left   = 195.76292338582678
$shp.Top    = 438.9848231496063
$shp.Width  = 347.1069191338583
$shp.Height = 73.51094488188977

# --- Change the label text from "first element" to "zero element" ---
# The original run text is:  "first element"  (curly quotes, 15 characters)
#   1      : opening curly quote
#   2-6    : "first"
#   7-14   : " element"
#   15     : closing curly quote
# Replace just the word "first" (characters 2-6) with "zero", which splits
# the text into three runs: the quote, "zero", and " element" + quote.
$tr = $shp.TextFrame.TextRange
$tr.Characters(2, 5).Text = "zero"
